# A new daily price record was inserted into the sheet at row 287, pushing
# every following record down by one row (287->288, 288->289, ... 380->381).
# The workbook's used range grows from A1:R380 to A1:R381 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 287; Excel shifts rows
# 287..380 down to 288..381 and the sheet dimension auto-extends to R381.
$ws.Rows.Item(287).Insert()

# Populate the newly inserted row 287 with the new record's data.
$ws.Range("A287").Value = 4
$ws.Range("B287").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C287").Value = "Los Lagos"
$ws.Range("D287").Value = 44876
$ws.Range("E287").Value = 10
$ws.Range("F287").Value = 100112045
$ws.Range("G287").Value = "Zapallo"
$ws.Range("H287").Value = "Paine"
$ws.Range("I287").Value = "1a (guarda)"
$ws.Range("J287").Value = 1000
$ws.Range("K287").Value = 700
$ws.Range("L287").Value = 750
$ws.Range("M287").Value = 725
$ws.Range("N287").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O287").Value = "Región de O'Higgins"
$ws.Range("P287").Value = 725
$ws.Range("Q287").Value = 1
$ws.Range("R287").Value = "Hortaliza"
